# Adds a new worksheet "Plan1" after the existing "Planilha1" sheet that
# holds a JXLS-style template cell ("${pessoas}") together with the
# explanatory author comment that JXLS/Excel authors typically leave on the
# template's anchor cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Create the new sheet right after the current (only) sheet so it ends up
# second in the workbook, matching sheetId/r:id ordering.
$plan1 = $wb.Worksheets.Add($null, $ws1)
$plan1.Name = "Plan1"

# Template placeholder cell.
$plan1.Range("A1").Value = '${pessoas}'

# Auto-fit column A to the new content, like Excel does after typing into it.
$plan1.Columns.Item(1).AutoFit()

# Use metric (cm based) page margins, as in the source workbook.
$plan1.PageSetup.LeftMargin = 36.850393728
$plan1.PageSetup.RightMargin = 36.850393728
$plan1.PageSetup.TopMargin = 56.692913399999995
$plan1.PageSetup.BottomMargin = 56.692913399999995
$plan1.PageSetup.HeaderMargin = 22.67716464
$plan1.PageSetup.FooterMargin = 22.67716464

# Author's explanatory comment describing the JXLS template directives.
$commentText = "Autor:`r`njx:area(lastCell=`"B2`")`r`njx:each(items=`"pessoas`" groupBy=`"name`" lastCell=`"B2`")`r`njx:each(items=`"_group.items`" var=`"pessoas`" lastCell=`"B1`")"
$comment = $plan1.Range("A1").AddComment($commentText)

# Leave the selection where the author left it.
$plan1.Range("K10").Select()
